$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 16824
$ws.Range("F2").Value = 0.662856216238388
$ws.Range("G2").Value = 4536.19637538216
$ws.Range("H2").Value = -9032.39275076432
$ws.Range("K2").Value = -8877.78151641667
$ws.Range("L2").Value = 7916.95827747725
$ws.Range("M2").Value = 16804

# Row 3
$ws.Range("E3").Value = 17133
$ws.Range("F3").Value = 0.668656623552693
$ws.Range("G3").Value = 4275.30267579662
$ws.Range("H3").Value = -8538.60535159325
$ws.Range("I3").Value = 493.787399171073
$ws.Range("J3").Value = [double]"5.96250358707656e-108"
$ws.Range("K3").Value = -8492.11278134925
$ws.Range("L3").Value = 8213.46993103883
$ws.Range("M3").Value = 17127

# Row 4
$ws.Range("E4").Value = 16824
$ws.Range("F4").Value = 0.792150335371764
$ws.Range("G4").Value = 1026.69043939897
$ws.Range("H4").Value = -2039.38087879795
$ws.Range("I4").Value = 6993.01187196637
$ws.Range("K4").Value = -1985.26694677627
$ws.Range("L4").Value = 11613.5269619364
$ws.Range("M4").Value = 16817
